$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) cells to remain plain text so values like
# "0.0770" / "63.095.46" are not auto-converted/rounded into numbers by Excel.
$priceCells = @("D2","D3","D5","D6","D7","D10","D11","D12","D13","D14","D16","D18","D19","D22","D23","D24","D26","D28","D34","D35","D36","D39","D40","D41","D42","D44","D45","D47","D48","D49","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Rows 34 and 35 swap places (EthereumClassic <-> Fetch.AI) with updated values
$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").Value = '1.34'
$ws.Range("E34").Value = '  -4.11%  '

$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").Value = '23.38'
$ws.Range("E35").Value = '  -0.65%  '

# Remaining per-cell Price / Volume(1h) updates
$ws.Range("D2").Value = '63.095.46'
$ws.Range("E2").Value = '  +1.93%  '
$ws.Range("D3").Value = '3.458.18'
$ws.Range("E3").Value = '  +1.26%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '580.22'
$ws.Range("E5").Value = '  +0.35%  '
$ws.Range("D6").Value = '147.54'
$ws.Range("E6").Value = '  +2.31%  '
$ws.Range("D7").Value = '3.458.10'
$ws.Range("E7").Value = '  +1.24%  '
$ws.Range("E9").Value = '  +0.60%  '
$ws.Range("D10").Value = '7.81'
$ws.Range("E10").Value = '  +2.86%  '
$ws.Range("D11").Value = '0.124'
$ws.Range("E11").Value = '  +1.08%  '
$ws.Range("D12").Value = '0.404'
$ws.Range("E12").Value = '  +4.89%  '
$ws.Range("D13").Value = '4.053.21'
$ws.Range("E13").Value = '  +1.29%  '
$ws.Range("D14").Value = '29.32'
$ws.Range("E14").Value = '  +3.27%  '
$ws.Range("E15").Value = '  +2.50%  '
$ws.Range("D16").Value = '3.472.92'
$ws.Range("E16").Value = '  +1.58%  '
$ws.Range("E17").Value = '  +1.09%  '
$ws.Range("D18").Value = '63.102.27'
$ws.Range("E18").Value = '  +1.88%  '
$ws.Range("D19").Value = '6.40'
$ws.Range("E19").Value = '  +4.00%  '
$ws.Range("E20").Value = '  +3.50%  '
$ws.Range("E21").Value = '  +1.56%  '
$ws.Range("D22").Value = '388.05'
$ws.Range("E22").Value = '  -0.51%  '
$ws.Range("D23").Value = '0.563'
$ws.Range("E23").Value = '  +1.70%  '
$ws.Range("D24").Value = '74.53'
$ws.Range("E24").Value = '  -0.52%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").Value = '3.608.20'
$ws.Range("E26").Value = '  +1.46%  '
$ws.Range("E27").Value = '  +1.08%  '
$ws.Range("D28").Value = '0.182'
$ws.Range("E28").Value = '  -1.33%  '
$ws.Range("E29").Value = '  +2.46%  '
$ws.Range("E30").Value = '  -0.10%  '
$ws.Range("E31").Value = '  +2.03%  '
$ws.Range("E32").Value = '  -0.16%  '
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("D36").Value = '5.32'
$ws.Range("E36").Value = '  +1.42%  '
$ws.Range("E37").Value = '  +2.59%  '
$ws.Range("E38").Value = '  +4.70%  '
$ws.Range("D39").Value = '31.87'
$ws.Range("E39").Value = '  +11.22%  '
$ws.Range("D40").Value = '168.27'
$ws.Range("E40").Value = '  +0.37%  '
$ws.Range("D41").Value = '3.495.60'
$ws.Range("E41").Value = '  +1.41%  '
$ws.Range("D42").Value = '0.0770'
$ws.Range("E42").Value = '  +2.14%  '
$ws.Range("E43").Value = '  +0.85%  '
$ws.Range("D44").Value = '1.73'
$ws.Range("E44").Value = '  +3.73%  '
$ws.Range("D45").Value = '42.38'
$ws.Range("E45").Value = '  -0.77%  '
$ws.Range("E46").Value = '  +3.58%  '
$ws.Range("D47").Value = '4.36'
$ws.Range("E47").Value = '  -1.30%  '
$ws.Range("D48").Value = '2.593.18'
$ws.Range("E48").Value = '  +3.60%  '
$ws.Range("D49").Value = '2.32'
$ws.Range("E49").Value = '  +11.59%  '
$ws.Range("D50").Value = '6.81'
$ws.Range("E50").Value = '  +2.75%  '
$ws.Range("D51").Value = '22.97'
$ws.Range("E51").Value = '  +0.85%  '
